$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Thu Jan 25 17:37:05 EST 2024"
$ws.Range("B3").Value = "Thu Jan 25 17:37:18 EST 2024"
$ws.Range("B4").Value = "Thu Jan 25 17:37:30 EST 2024"
$ws.Range("B5").Value = "Thu Jan 25 17:37:43 EST 2024"
$ws.Range("B6").Value = "Thu Jan 25 17:37:55 EST 2024"
$ws.Range("B7").Value = "Thu Jan 25 17:38:07 EST 2024"
$ws.Range("B8").Value = "Thu Jan 25 17:38:20 EST 2024"
$ws.Range("B9").Value = "Thu Jan 25 17:38:32 EST 2024"
$ws.Range("B10").Value = "Thu Jan 25 17:38:44 EST 2024"
$ws.Range("B11").Value = "Thu Jan 25 17:38:57 EST 2024"
$ws.Range("B12").Value = "Thu Jan 25 17:39:09 EST 2024"
$ws.Range("B13").Value = "Thu Jan 25 17:39:21 EST 2024"
$ws.Range("B14").Value = "Thu Jan 25 17:39:34 EST 2024"
$ws.Range("B15").Value = "Thu Jan 25 17:39:46 EST 2024"
$ws.Range("B16").Value = "Thu Jan 25 17:39:59 EST 2024"
$ws.Range("B17").Value = "Thu Jan 25 17:40:11 EST 2024"
$ws.Range("B18").Value = "Thu Jan 25 17:40:24 EST 2024"
$ws.Range("B19").Value = "Thu Jan 25 17:40:36 EST 2024"
$ws.Range("B20").Value = "Thu Jan 25 17:40:48 EST 2024"
$ws.Range("B28").Value = "Thu Jan 25 17:41:02 EST 2024"
$ws.Range("B29").Value = "Thu Jan 25 17:41:14 EST 2024"
$ws.Range("B30").Value = "Thu Jan 25 17:41:26 EST 2024"
$ws.Range("B31").Value = "Thu Jan 25 17:41:39 EST 2024"
$ws.Range("B32").Value = "Thu Jan 25 17:41:51 EST 2024"
$ws.Range("B33").Value = "Thu Jan 25 17:42:03 EST 2024"
$ws.Range("B34").Value = "Thu Jan 25 17:42:15 EST 2024"
$ws.Range("B35").Value = "Thu Jan 25 17:42:27 EST 2024"
$ws.Range("B36").Value = "Thu Jan 25 17:42:40 EST 2024"
$ws.Range("B37").Value = "Thu Jan 25 17:42:55 EST 2024"
$ws.Range("B38").Value = "Thu Jan 25 17:43:07 EST 2024"
$ws.Range("B39").Value = "Thu Jan 25 17:43:19 EST 2024"
$ws.Range("B40").Value = "Thu Jan 25 17:43:32 EST 2024"
$ws.Range("B41").Value = "Thu Jan 25 17:43:44 EST 2024"
$ws.Range("B42").Value = "Thu Jan 25 17:43:57 EST 2024"
$ws.Range("B43").Value = "Thu Jan 25 17:44:09 EST 2024"
$ws.Range("B44").Value = "Thu Jan 25 17:44:21 EST 2024"
$ws.Range("B45").Value = "Thu Jan 25 17:44:34 EST 2024"
$ws.Range("B46").Value = "Thu Jan 25 17:44:46 EST 2024"
$ws.Range("B47").Value = "Thu Jan 25 17:45:00 EST 2024"
$ws.Range("B48").Value = "Thu Jan 25 17:45:12 EST 2024"
$ws.Range("B49").Value = "Thu Jan 25 17:45:25 EST 2024"
$ws.Range("B50").Value = "Thu Jan 25 17:45:37 EST 2024"
$ws.Range("B51").Value = "Thu Jan 25 17:45:49 EST 2024"
$ws.Range("B52").Value = "Thu Jan 25 17:46:01 EST 2024"
$ws.Range("B53").Value = "Thu Jan 25 17:46:14 EST 2024"
$ws.Range("B54").Value = "Thu Jan 25 17:46:26 EST 2024"
